$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.4800474766800295
$ws.Range("J2").Value = 0.4800474766800294
$ws.Range("M2").Value = 0.6946430000000001
$ws.Range("N2").Value = 2.083929
$ws.Range("O2").Value = 0.1140293552421611
$ws.Range("P2").Value = 0.1140293552421611
$ws.Range("Q2").Value = 0.4881934795663334
$ws.Range("R2").Value = 4.393741316097
$ws.Range("S2").Value = 0.05473950425145014
$ws.Range("T2").Value = 0.05473950425145013

# Row 3
$ws.Range("I3").Value = 0.4800474766800295
$ws.Range("J3").Value = 0.4800474766800294
$ws.Range("M3").Value = 5.373609333333333
$ws.Range("N3").Value = 16.120828
$ws.Range("O3").Value = 0.8821066470161785
$ws.Range("P3").Value = 0.8821066470161785
$ws.Range("Q3").Value = 3.776560101044889
$ws.Range("R3").Value = 33.989040909404
$ws.Range("S3").Value = 0.4234530700627979
$ws.Range("T3").Value = 0.4234530700627979

# Row 4
$ws.Range("I4").Value = 0.4800474766800295
$ws.Range("J4").Value = 0.4800474766800294
$ws.Range("M4").Value = 0.02353866666666667
$ws.Range("N4").Value = 0.070616
$ws.Range("O4").Value = 0.00386399774166032
$ws.Range("P4").Value = 0.00386399774166032
$ws.Range("Q4").Value = 0.01654292000977778
$ws.Range("R4").Value = 0.148886280088
$ws.Range("S4").Value = 0.001854902365781369
$ws.Range("T4").Value = 0.001854902365781369

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.341519
$ws.Range("H5").Value = 1.024557
$ws.Range("I5").Value = 0.233275296666637
$ws.Range("J5").Value = 0.233275296666637
$ws.Range("M5").Value = 0.6946430000000001
$ws.Range("N5").Value = 2.083929
$ws.Range("O5").Value = 0.1140293552421611
$ws.Range("P5").Value = 0.1140293552421611
$ws.Range("Q5").Value = 0.237233782717
$ws.Range("R5").Value = 2.135104044453
$ws.Range("S5").Value = 0.02660023167282048
$ws.Range("T5").Value = 0.02660023167282048

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.341519
$ws.Range("H6").Value = 1.024557
$ws.Range("I6").Value = 0.233275296666637
$ws.Range("J6").Value = 0.233275296666637
$ws.Range("M6").Value = 5.373609333333333
$ws.Range("N6").Value = 16.120828
$ws.Range("O6").Value = 0.8821066470161785
$ws.Range("P6").Value = 0.8821066470161785
$ws.Range("Q6").Value = 1.835189685910666
$ws.Range("R6").Value = 16.516707173196
$ws.Range("S6").Value = 0.2057736897743115
$ws.Range("T6").Value = 0.2057736897743115

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.341519
$ws.Range("H7").Value = 1.024557
$ws.Range("I7").Value = 0.233275296666637
$ws.Range("J7").Value = 0.233275296666637
$ws.Range("M7").Value = 0.02353866666666667
$ws.Range("N7").Value = 0.070616
$ws.Range("O7").Value = 0.00386399774166032
$ws.Range("P7").Value = 0.00386399774166032
$ws.Range("Q7").Value = 0.008038901901333332
$ws.Range("R7").Value = 0.072350117112
$ws.Range("S7").Value = 0.0009013752195050267
$ws.Range("T7").Value = 0.0009013752195050267

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4197003333333333
$ws.Range("H8").Value = 1.259101
$ws.Range("I8").Value = 0.2866772266533335
$ws.Range("J8").Value = 0.2866772266533335
$ws.Range("M8").Value = 0.6946430000000001
$ws.Range("N8").Value = 2.083929
$ws.Range("O8").Value = 0.1140293552421611
$ws.Range("P8").Value = 0.1140293552421611
$ws.Range("Q8").Value = 0.2915418986476667
$ws.Range("R8").Value = 2.623877087829
$ws.Range("S8").Value = 0.03268961931789051
$ws.Range("T8").Value = 0.0326896193178905

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4197003333333333
$ws.Range("H9").Value = 1.259101
$ws.Range("I9").Value = 0.2866772266533335
$ws.Range("J9").Value = 0.2866772266533335
$ws.Range("M9").Value = 5.373609333333333
$ws.Range("N9").Value = 16.120828
$ws.Range("O9").Value = 0.8821066470161785
$ws.Range("P9").Value = 0.8821066470161785
$ws.Range("Q9").Value = 2.255305628403111
$ws.Range("R9").Value = 20.297750655628
$ws.Range("S9").Value = 0.2528798871790691
$ws.Range("T9").Value = 0.252879887179069

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4197003333333333
$ws.Range("H10").Value = 1.259101
$ws.Range("I10").Value = 0.2866772266533335
$ws.Range("J10").Value = 0.2866772266533335
$ws.Range("M10").Value = 0.02353866666666667
$ws.Range("N10").Value = 0.070616
$ws.Range("O10").Value = 0.00386399774166032
$ws.Range("P10").Value = 0.00386399774166032
$ws.Range("Q10").Value = 0.009879186246222223
$ws.Range("R10").Value = 0.088912676216
$ws.Range("S10").Value = 0.001107720156373924
$ws.Range("T10").Value = 0.001107720156373924
